$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 1762
    8  = 103
    11 = 1457
    13 = 577
    14 = 361
    20 = 130
    23 = 3436
    24 = 413
    25 = 298
    27 = 74
    28 = 22
    30 = 1168
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
